$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 58

# Columns whose values look numeric/date-like need an explicit text marker
# (leading apostrophe) so Excel stores them as text, matching the source
# data export (Caso/OT/Comuna codes and the US-style date string) instead of
# auto-converting them to numbers or a date serial.
$ws.Cells.Item($row, 1).Value  = "'6484"
$ws.Cells.Item($row, 2).Value  = "'7/24/2025"
$ws.Cells.Item($row, 3).Value  = "URIARTE 2396"
$ws.Cells.Item($row, 4).Value  = "'14"
$ws.Cells.Item($row, 5).Value  = "'808509373"
$ws.Cells.Item($row, 6).Value  = "NEW"
$ws.Cells.Item($row, 7).Value  = "Pendiente"
$ws.Cells.Item($row, 8).Value  = "Picada"
$ws.Cells.Item($row, 9).Value  = 1
$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Pasante"
$ws.Cells.Item($row, 13).Value = -58.423934
$ws.Cells.Item($row, 14).Value = -34.581562
$ws.Cells.Item($row, 15).Value = "Palermo"
$ws.Cells.Item($row, 16).Value = "Capital Sur"
